$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column G (Affichage mail / ${com.temShowMail}) -- this shifts H..M left by one
$ws.Columns("G").Delete()

# Insert a new column before column L (which currently holds "Alerte proposition"),
# so the old L (Adresse) stays in K and the new column becomes L: "Mail alerte" / ${com.mailAlert}
$ws.Columns("L").Insert()

# Match new column L's width/style to its neighbour K
$ws.Columns("L").ColumnWidth = $ws.Columns("K").ColumnWidth
$ws.Columns("L").Style = $ws.Columns("K").Style

$ws.Cells.Item(1, 12).Value = "Mail alerte"
$ws.Cells.Item(2, 12).Value = "`${com.mailAlert}"

for ($col = 1; $col -le 13; $col++) {
    Write-Output ("col=" + $col + " width=" + $ws.Columns($col).ColumnWidth)
}
